$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 8
$ws.Cells.Item(8, 8).Value = 1515
$ws.Cells.Item(8, 9).Value = 224
$ws.Cells.Item(8, 11).Value = 672
$ws.Cells.Item(8, 13).Value = -533
# row 10
$ws.Cells.Item(10, 8).Value = 4624.875
$ws.Cells.Item(10, 10).Value = 4999.857
$ws.Cells.Item(10, 12).Value = 4999.857
$ws.Cells.Item(10, 14).Value = -5585.857
# row 31
$ws.Cells.Item(31, 8).Value = 499.33334
$ws.Cells.Item(31, 9).Value = 499.33334
$ws.Cells.Item(31, 11).Value = 1498.00002
$ws.Cells.Item(31, 13).Value = -1268.00002
# row 48
$ws.Cells.Item(48, 8).Value = 1750
$ws.Cells.Item(48, 9).Value = 3000
$ws.Cells.Item(48, 10).Value = 500
$ws.Cells.Item(48, 11).Value = 9000
$ws.Cells.Item(48, 12).Value = 1500
$ws.Cells.Item(48, 13).Value = -8708
$ws.Cells.Item(48, 14).Value = -2084
# row 56
$ws.Cells.Item(56, 8).Value = 1750
$ws.Cells.Item(56, 9).Value = 3000
$ws.Cells.Item(56, 10).Value = 500
$ws.Cells.Item(56, 11).Value = 9000
$ws.Cells.Item(56, 12).Value = 1500
$ws.Cells.Item(56, 13).Value = -8466
$ws.Cells.Item(56, 14).Value = -2568
# row 62
$ws.Cells.Item(62, 8).Value = 1707.8572
$ws.Cells.Item(62, 9).Value = 1019
$ws.Cells.Item(62, 10).Value = 2090.5557
$ws.Cells.Item(62, 11).Value = 1019
$ws.Cells.Item(62, 12).Value = 2090.5557
$ws.Cells.Item(62, 13).Value = -395
$ws.Cells.Item(62, 14).Value = -3338.5557
# row 65
$ws.Cells.Item(65, 8).Value = 1707.8572
$ws.Cells.Item(65, 9).Value = 1019
$ws.Cells.Item(65, 10).Value = 2090.5557
$ws.Cells.Item(65, 11).Value = 5095
$ws.Cells.Item(65, 12).Value = 10452.7785
$ws.Cells.Item(65, 13).Value = -1975
$ws.Cells.Item(65, 14).Value = -16692.7785
# row 98
$ws.Cells.Item(98, 8).Value = 2493.4546
$ws.Cells.Item(98, 9).Value = 2652.5
$ws.Cells.Item(98, 10).Value = 903
$ws.Cells.Item(98, 11).Value = 2652.5
$ws.Cells.Item(98, 12).Value = 903
$ws.Cells.Item(98, 13).Value = -1154.5
$ws.Cells.Item(98, 14).Value = -3899
# row 116
$ws.Cells.Item(116, 8).Value = 3073.3333
$ws.Cells.Item(116, 9).Value = 3406.25
$ws.Cells.Item(116, 11).Value = 3406.25
$ws.Cells.Item(116, 13).Value = 35.75
# row 122
$ws.Cells.Item(122, 8).Value = 2493.4546
$ws.Cells.Item(122, 9).Value = 2652.5
$ws.Cells.Item(122, 10).Value = 903
$ws.Cells.Item(122, 11).Value = 7957.5
$ws.Cells.Item(122, 12).Value = 2709
$ws.Cells.Item(122, 13).Value = -5507.5
$ws.Cells.Item(122, 14).Value = -7609
# row 129
$ws.Cells.Item(129, 8).Value = 1130.6666
$ws.Cells.Item(129, 9).Value = 462
$ws.Cells.Item(129, 10).Value = 1465
$ws.Cells.Item(129, 11).Value = 1386
$ws.Cells.Item(129, 12).Value = 4395
$ws.Cells.Item(129, 13).Value = 3614
$ws.Cells.Item(129, 14).Value = -14395
# row 132
$ws.Cells.Item(132, 8).Value = 165804.69
$ws.Cells.Item(132, 9).Value = 168553.1
$ws.Cells.Item(132, 11).Value = 505659.3
$ws.Cells.Item(132, 13).Value = -503129.3
# row 137
$ws.Cells.Item(137, 8).Value = 6316.591
$ws.Cells.Item(137, 9).Value = 7779
$ws.Cells.Item(137, 10).Value = 2416.8333
$ws.Cells.Item(137, 11).Value = 23337
$ws.Cells.Item(137, 12).Value = 7250.499899999999
$ws.Cells.Item(137, 13).Value = -20787
$ws.Cells.Item(137, 14).Value = -12350.4999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Cells.Item(2, 8).Value = 3270572.2
$ws.Cells.Item(2, 9).Value = 2923.25
$ws.Cells.Item(2, 10).Value = 29411764
$ws.Cells.Item(2, 11).Value = 2923.25
$ws.Cells.Item(2, 12).Value = 29411764
$ws.Cells.Item(2, 13).Value = -2810.25
$ws.Cells.Item(2, 14).Value = -29411990
# row 4
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = $null
$ws.Cells.Item(4, 14).Value = $null
# row 6
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 14).Value = $null
# row 36
$ws.Cells.Item(36, 8).Value = 3666.6667
$ws.Cells.Item(36, 9).Value = 3666.6667
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 3666.6667
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -3320.6667
$ws.Cells.Item(36, 14).Value = $null
# row 41
$ws.Cells.Item(41, 8).Value = 15317.223
$ws.Cells.Item(41, 9).Value = 1285
$ws.Cells.Item(41, 10).Value = 22333.334
$ws.Cells.Item(41, 11).Value = 1285
$ws.Cells.Item(41, 12).Value = 22333.334
$ws.Cells.Item(41, 13).Value = -871
$ws.Cells.Item(41, 14).Value = -23161.334
# row 116
$ws.Cells.Item(116, 8).Value = 3270572.2
$ws.Cells.Item(116, 9).Value = 2923.25
$ws.Cells.Item(116, 10).Value = 29411764
$ws.Cells.Item(116, 11).Value = 2923.25
$ws.Cells.Item(116, 12).Value = 29411764
$ws.Cells.Item(116, 13).Value = -629.25
$ws.Cells.Item(116, 14).Value = -29416352
# row 122
$ws.Cells.Item(122, 8).Value = 12822360
$ws.Cells.Item(122, 9).Value = 17242940
$ws.Cells.Item(122, 10).Value = 2679.9
$ws.Cells.Item(122, 11).Value = 51728820
$ws.Cells.Item(122, 12).Value = 8039.700000000001
$ws.Cells.Item(122, 13).Value = -51726370
$ws.Cells.Item(122, 14).Value = -12939.7

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Cells.Item(3, 8).Value = 3270572.2
$ws.Cells.Item(3, 9).Value = 2923.25
$ws.Cells.Item(3, 10).Value = 29411764
$ws.Cells.Item(3, 11).Value = 2923.25
$ws.Cells.Item(3, 12).Value = 29411764
$ws.Cells.Item(3, 13).Value = -2809.25
$ws.Cells.Item(3, 14).Value = -29411992
# row 7
$ws.Cells.Item(7, 8).Value = 1670000
$ws.Cells.Item(7, 9).Value = 5000000
$ws.Cells.Item(7, 10).Value = 5000
$ws.Cells.Item(7, 11).Value = 5000000
$ws.Cells.Item(7, 12).Value = 5000
$ws.Cells.Item(7, 13).Value = -4999887
$ws.Cells.Item(7, 14).Value = -5226
# row 107
$ws.Cells.Item(107, 8).Value = 1344.1052
$ws.Cells.Item(107, 9).Value = 1068.7142
$ws.Cells.Item(107, 10).Value = 2115.2
$ws.Cells.Item(107, 11).Value = 1068.7142
$ws.Cells.Item(107, 12).Value = 2115.2
$ws.Cells.Item(107, 13).Value = 851.2858000000001
$ws.Cells.Item(107, 14).Value = -5955.2
# row 135
$ws.Cells.Item(135, 8).Value = 57463.625
$ws.Cells.Item(135, 10).Value = 57463.625
$ws.Cells.Item(135, 12).Value = 57463.625
$ws.Cells.Item(135, 14).Value = -67603.625

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Cells.Item(7, 8).Value = 739
$ws.Cells.Item(7, 9).Value = 424.375
$ws.Cells.Item(7, 10).Value = 1997.5
$ws.Cells.Item(7, 11).Value = 424.375
$ws.Cells.Item(7, 12).Value = 1997.5
$ws.Cells.Item(7, 13).Value = -311.375
$ws.Cells.Item(7, 14).Value = -2223.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 4
$ws.Cells.Item(4, 8).Value = 625.875
$ws.Cells.Item(4, 9).Value = 280.07144
$ws.Cells.Item(4, 11).Value = 840.21432
$ws.Cells.Item(4, 13).Value = -728.21432
# row 5
$ws.Cells.Item(5, 8).Value = 385.30304
$ws.Cells.Item(5, 9).Value = 382.875
$ws.Cells.Item(5, 11).Value = 1148.625
$ws.Cells.Item(5, 13).Value = -1036.625
# row 68
$ws.Cells.Item(68, 8).Value = 1006.9452
$ws.Cells.Item(68, 9).Value = 761.4666999999999
$ws.Cells.Item(68, 10).Value = 1070.431
$ws.Cells.Item(68, 11).Value = 2284.4001
$ws.Cells.Item(68, 12).Value = 3211.293
$ws.Cells.Item(68, 13).Value = -1473.4001
$ws.Cells.Item(68, 14).Value = -4833.293
# row 71
$ws.Cells.Item(71, 8).Value = 1006.9452
$ws.Cells.Item(71, 9).Value = 761.4666999999999
$ws.Cells.Item(71, 10).Value = 1070.431
$ws.Cells.Item(71, 11).Value = 6853.2003
$ws.Cells.Item(71, 12).Value = 9633.879000000001
$ws.Cells.Item(71, 13).Value = -2797.2003
$ws.Cells.Item(71, 14).Value = -17745.879
# row 101
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(101, 14).Value = $null
# row 135
$ws.Cells.Item(135, 8).Value = 385.30304
$ws.Cells.Item(135, 9).Value = 382.875
$ws.Cells.Item(135, 11).Value = 3445.875
$ws.Cells.Item(135, 13).Value = -910.875

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 39
$ws.Cells.Item(39, 8).Value = 32500
$ws.Cells.Item(39, 10).Value = 32500
$ws.Cells.Item(39, 12).Value = 32500
$ws.Cells.Item(39, 14).Value = -33564

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 21
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).Value = $null
# row 22
$ws.Cells.Item(22, 8).Value = 283.73685
$ws.Cells.Item(22, 9).Value = 320.53845
$ws.Cells.Item(22, 10).Value = 204
$ws.Cells.Item(22, 11).Value = 320.53845
$ws.Cells.Item(22, 12).Value = 204
$ws.Cells.Item(22, 13).Value = -25.53845000000001
$ws.Cells.Item(22, 14).Value = -794
# row 27
$ws.Cells.Item(27, 8).Value = 283.73685
$ws.Cells.Item(27, 9).Value = 320.53845
$ws.Cells.Item(27, 10).Value = 204
$ws.Cells.Item(27, 11).Value = 320.53845
$ws.Cells.Item(27, 12).Value = 204
$ws.Cells.Item(27, 13).Value = -213.53845
$ws.Cells.Item(27, 14).Value = -418
# row 46
$ws.Cells.Item(46, 8).Value = 724.8421
$ws.Cells.Item(46, 9).Value = 701.6923
$ws.Cells.Item(46, 11).Value = 701.6923
$ws.Cells.Item(46, 13).Value = -513.6923
# row 61
$ws.Cells.Item(61, 8).Value = 3312.25
$ws.Cells.Item(61, 9).Value = 2374.75
$ws.Cells.Item(61, 10).Value = 4249.75
$ws.Cells.Item(61, 11).Value = 2374.75
$ws.Cells.Item(61, 12).Value = 4249.75
$ws.Cells.Item(61, 13).Value = -2172.75
$ws.Cells.Item(61, 14).Value = -4653.75
# row 113
$ws.Cells.Item(113, 8).Value = 3312.25
$ws.Cells.Item(113, 9).Value = 2374.75
$ws.Cells.Item(113, 10).Value = 4249.75
$ws.Cells.Item(113, 11).Value = 2374.75
$ws.Cells.Item(113, 12).Value = 4249.75
$ws.Cells.Item(113, 13).Value = -204.75
$ws.Cells.Item(113, 14).Value = -8589.75

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 40
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = $null
$ws.Cells.Item(40, 14).Value = $null
# row 107
$ws.Cells.Item(107, 8).Value = 1370.3334
$ws.Cells.Item(107, 9).Value = 1061.5385
$ws.Cells.Item(107, 10).Value = 1657.0714
$ws.Cells.Item(107, 11).Value = 3184.6155
$ws.Cells.Item(107, 12).Value = 4971.2142
$ws.Cells.Item(107, 13).Value = -1264.6155
$ws.Cells.Item(107, 14).Value = -8811.2142
# row 132
$ws.Cells.Item(132, 8).Value = 5020.387
$ws.Cells.Item(132, 9).Value = 5896.913
$ws.Cells.Item(132, 10).Value = 2500.375
$ws.Cells.Item(132, 11).Value = 17690.739
$ws.Cells.Item(132, 12).Value = 7501.125
$ws.Cells.Item(132, 13).Value = -15160.739
$ws.Cells.Item(132, 14).Value = -12561.125
